$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 9 data: E9 and F9 go from 0 to 1 (I9's formula recalculates)
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1

# Add new row 10 of data (all zeros) with the extended formula in I10
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Formula = '=A10*$A$1+B10*$B$1+C10*$C$1+D10*$D$1+E10*$E$1+F10*$F$1+G10*$G$1+H10*$H$1'

# New conditional formatting rule for the new row A10:H10 (same rule as the existing one: highlight cells equal to 1)
$newRange = $ws.Range("A10:H10")
$newRule = $newRange.FormatConditions.Add(1, 3, "1")
$newRule.Interior.ColorIndex = 1
$newRule.SetFirstPriority()

# Move the selection to E10 to match the author's final cursor position
$ws.Range("E10").Select() | Out-Null
